$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valA = $cellA.Value2
        $valB = $cellB.Value2

        $cellA.Value2 = $valB
        $cellB.Value2 = $valA
    }
}

# Swap data (columns B=2 through AB=28) between rows 11 and 12, keeping column A untouched
Swap-Rows 11 12 2 28

# Swap data (columns B=2 through AB=28) between rows 83 and 84, keeping column A untouched
Swap-Rows 83 84 2 28
